# This script reproduces the commit "added example and template files":
# the AusDiab example/data rows are stripped from every sheet so that only
# the field-name header row (row 1) and the field-description row (row 2)
# remain, turning the filled-in example workbook back into a blank
# template.

$wb = $excel.ActiveWorkbook

# --- subject sheet -------------------------------------------------------
# The "subject" sheet already keeps a couple of blank, pre-styled template
# rows below the example data (rows 5-6), so here the example values are
# simply cleared out (not deleted) to match that same blank-template look.
$wsSubject = $wb.Worksheets.Item("subject")
$wsSubject.Range("A3:D4").ClearContents()
$wsSubject.Range("A3:E4").Select()

# --- sample sheet ---------------------------------------------------------
$wsSample = $wb.Worksheets.Item("sample")
$wsSample.Rows("3:5").Delete()
$wsSample.Range("A3:N5").Select()

# --- lipidomics_assay sheet ------------------------------------------------
$wsAssay = $wb.Worksheets.Item("lipidomics_assay")
$wsAssay.Rows("3:5").Delete()
$wsAssay.Range("A3:E5").Select()

# --- lipidomics_file sheet -------------------------------------------------
$wsFile = $wb.Worksheets.Item("lipidomics_file")
$wsFile.Rows("3:4").Delete()
$wsFile.Range("A3:L4").Select()

# --- lipidomics_mapping_file sheet -----------------------------------------
$wsMap = $wb.Worksheets.Item("lipidomics_mapping_file")
$wsMap.Rows("3:4").Delete()
$wsMap.Activate()
$wsMap.Application.ActiveWindow.ScrollColumn = 1
$wsMap.Range("B11").Select()
